# Apply full data scrape update for extra batting and bowling fields
# Target sheet: "ODI Batting Extra"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Each entry: row, column letter, kind (TXT/NUM/EMPTY), value, protect-as-text flag
$data = @(
    @(2, 'A', 'TXT', '3393', 1),
    @(2, 'B', 'NUM', '11', 0),
    @(2, 'C', 'TXT', '0', 1),
    @(2, 'D', 'TXT', '0', 1),
    @(2, 'E', 'EMPTY', '', 0),
    @(2, 'F', 'TXT', 'NO', 0),
    @(3, 'A', 'TXT', '3402', 1),
    @(3, 'B', 'NUM', '11', 0),
    @(3, 'C', 'TXT', '1', 1),
    @(3, 'D', 'TXT', '0', 1),
    @(3, 'E', 'TXT', '1.82%', 1),
    @(3, 'F', 'TXT', 'NO', 0),
    @(4, 'A', 'TXT', '3659', 1),
    @(4, 'B', 'NUM', '11', 0),
    @(4, 'C', 'EMPTY', '', 0),
    @(4, 'D', 'EMPTY', '', 0),
    @(4, 'E', 'EMPTY', '', 0),
    @(4, 'F', 'TXT', 'NO', 0),
    @(5, 'A', 'TXT', '3668', 1),
    @(5, 'B', 'NUM', '11', 0),
    @(5, 'C', 'TXT', '0', 1),
    @(5, 'D', 'TXT', '1', 1),
    @(5, 'E', 'TXT', '3.83%', 1),
    @(5, 'F', 'TXT', 'NO', 0),
    @(6, 'A', 'TXT', '3669', 1),
    @(6, 'B', 'NUM', '11', 0),
    @(6, 'C', 'EMPTY', '', 0),
    @(6, 'D', 'EMPTY', '', 0),
    @(6, 'E', 'EMPTY', '', 0),
    @(6, 'F', 'TXT', 'NO', 0),
    @(7, 'A', 'TXT', '3673', 1),
    @(7, 'B', 'EMPTY', '', 0),
    @(7, 'C', 'EMPTY', '', 0),
    @(7, 'D', 'EMPTY', '', 0),
    @(7, 'E', 'EMPTY', '', 0),
    @(7, 'F', 'TXT', 'NO', 0),
    @(8, 'A', 'TXT', '3677', 1),
    @(8, 'B', 'EMPTY', '', 0),
    @(8, 'C', 'EMPTY', '', 0),
    @(8, 'D', 'EMPTY', '', 0),
    @(8, 'E', 'EMPTY', '', 0),
    @(8, 'F', 'TXT', 'NO', 0),
    @(9, 'A', 'TXT', '3679', 1),
    @(9, 'B', 'EMPTY', '', 0),
    @(9, 'C', 'EMPTY', '', 0),
    @(9, 'D', 'EMPTY', '', 0),
    @(9, 'E', 'EMPTY', '', 0),
    @(9, 'F', 'TXT', 'NO', 0),
    @(10, 'A', 'TXT', '3877', 1),
    @(10, 'B', 'EMPTY', '', 0),
    @(10, 'C', 'EMPTY', '', 0),
    @(10, 'D', 'EMPTY', '', 0),
    @(10, 'E', 'EMPTY', '', 0),
    @(10, 'F', 'TXT', 'NO', 0),
    @(11, 'A', 'TXT', '3878', 1),
    @(11, 'B', 'EMPTY', '', 0),
    @(11, 'C', 'EMPTY', '', 0),
    @(11, 'D', 'EMPTY', '', 0),
    @(11, 'E', 'EMPTY', '', 0),
    @(11, 'F', 'TXT', 'NO', 0),
    @(12, 'A', 'TXT', '3893', 1),
    @(12, 'B', 'NUM', '10', 0),
    @(12, 'C', 'TXT', '4', 1),
    @(12, 'D', 'TXT', '0', 1),
    @(12, 'E', 'TXT', '21.13%', 1),
    @(12, 'F', 'TXT', 'NO', 0),
    @(13, 'A', 'TXT', '3894', 1),
    @(13, 'B', 'EMPTY', '', 0),
    @(13, 'C', 'EMPTY', '', 0),
    @(13, 'D', 'EMPTY', '', 0),
    @(13, 'E', 'EMPTY', '', 0),
    @(13, 'F', 'TXT', 'NO', 0),
    @(14, 'A', 'TXT', '3924', 1),
    @(14, 'B', 'NUM', '10', 0),
    @(14, 'C', 'TXT', '0', 1),
    @(14, 'D', 'TXT', '0', 1),
    @(14, 'E', 'TXT', '0.49%', 1),
    @(14, 'F', 'TXT', 'NO', 0),
    @(15, 'A', 'TXT', '4169', 1),
    @(15, 'B', 'NUM', '9', 0),
    @(15, 'C', 'EMPTY', '', 0),
    @(15, 'D', 'EMPTY', '', 0),
    @(15, 'E', 'EMPTY', '', 0),
    @(15, 'F', 'TXT', 'NO', 0),
    @(16, 'A', 'TXT', '4170', 1),
    @(16, 'B', 'EMPTY', '', 0),
    @(16, 'C', 'EMPTY', '', 0),
    @(16, 'D', 'EMPTY', '', 0),
    @(16, 'E', 'EMPTY', '', 0),
    @(16, 'F', 'TXT', 'NO', 0),
    @(17, 'A', 'TXT', '4234', 1),
    @(17, 'B', 'NUM', '10', 0),
    @(17, 'C', 'TXT', '1', 1),
    @(17, 'D', 'TXT', '0', 1),
    @(17, 'E', 'TXT', '2.48%', 1),
    @(17, 'F', 'TXT', 'NO', 0),
    @(18, 'A', 'TXT', '4235', 1),
    @(18, 'B', 'NUM', '10', 0),
    @(18, 'C', 'EMPTY', '', 0),
    @(18, 'D', 'EMPTY', '', 0),
    @(18, 'E', 'EMPTY', '', 0),
    @(18, 'F', 'TXT', 'NO', 0),
    @(19, 'A', 'TXT', '4263', 1),
    @(19, 'B', 'NUM', '10', 0),
    @(19, 'C', 'TXT', '0', 1),
    @(19, 'D', 'TXT', '0', 1),
    @(19, 'E', 'TXT', '0.37%', 1),
    @(19, 'F', 'TXT', 'NO', 0),
    @(20, 'A', 'TXT', '4266', 1),
    @(20, 'B', 'NUM', '9', 0),
    @(20, 'C', 'EMPTY', '', 0),
    @(20, 'D', 'EMPTY', '', 0),
    @(20, 'E', 'EMPTY', '', 0),
    @(20, 'F', 'TXT', 'NO', 0),
    @(21, 'A', 'TXT', '4270', 1),
    @(21, 'B', 'NUM', '9', 0),
    @(21, 'C', 'EMPTY', '', 0),
    @(21, 'D', 'EMPTY', '', 0),
    @(21, 'E', 'EMPTY', '', 0),
    @(21, 'F', 'TXT', 'NO', 0),
    @(22, 'A', 'TXT', '4273', 1),
    @(22, 'B', 'NUM', '9', 0),
    @(22, 'C', 'EMPTY', '', 0),
    @(22, 'D', 'EMPTY', '', 0),
    @(22, 'E', 'EMPTY', '', 0),
    @(22, 'F', 'TXT', 'NO', 0),
    @(23, 'A', 'TXT', '4274', 1),
    @(23, 'B', 'EMPTY', '', 0),
    @(23, 'C', 'EMPTY', '', 0),
    @(23, 'D', 'EMPTY', '', 0),
    @(23, 'E', 'EMPTY', '', 0),
    @(23, 'F', 'TXT', 'NO', 0),
    @(24, 'A', 'TXT', '4275', 1),
    @(24, 'B', 'NUM', '10', 0),
    @(24, 'C', 'EMPTY', '', 0),
    @(24, 'D', 'EMPTY', '', 0),
    @(24, 'E', 'EMPTY', '', 0),
    @(24, 'F', 'TXT', 'NO', 0),
    @(25, 'A', 'TXT', '4276', 1),
    @(25, 'B', 'NUM', '10', 0),
    @(25, 'C', 'EMPTY', '', 0),
    @(25, 'D', 'EMPTY', '', 0),
    @(25, 'E', 'EMPTY', '', 0),
    @(25, 'F', 'TXT', 'NO', 0),
    @(26, 'A', 'TXT', '4277', 1),
    @(26, 'B', 'NUM', '11', 0),
    @(26, 'C', 'TXT', '0', 1),
    @(26, 'D', 'TXT', '0', 1),
    @(26, 'E', 'EMPTY', '', 0),
    @(26, 'F', 'TXT', 'NO', 0),
    @(27, 'A', 'TXT', '4336', 1),
    @(27, 'B', 'NUM', '11', 0),
    @(27, 'C', 'TXT', '0', 1),
    @(27, 'D', 'TXT', '0', 1),
    @(27, 'E', 'TXT', '2.24%', 1),
    @(27, 'F', 'TXT', 'NO', 0),
    @(28, 'A', 'TXT', '4341', 1),
    @(28, 'B', 'EMPTY', '', 0),
    @(28, 'C', 'EMPTY', '', 0),
    @(28, 'D', 'EMPTY', '', 0),
    @(28, 'E', 'EMPTY', '', 0),
    @(28, 'F', 'EMPTY', '', 0),
    @(29, 'A', 'TXT', '4351', 1),
    @(29, 'B', 'EMPTY', '', 0),
    @(29, 'C', 'EMPTY', '', 0),
    @(29, 'D', 'EMPTY', '', 0),
    @(29, 'E', 'EMPTY', '', 0),
    @(29, 'F', 'EMPTY', '', 0),
    @(30, 'A', 'TXT', '4354', 1),
    @(30, 'B', 'EMPTY', '', 0),
    @(30, 'C', 'EMPTY', '', 0),
    @(30, 'D', 'EMPTY', '', 0),
    @(30, 'E', 'EMPTY', '', 0),
    @(30, 'F', 'EMPTY', '', 0)
)

foreach ($item in $data) {
    $r = $item[0]
    $col = $item[1]
    $kind = $item[2]
    $val = $item[3]
    $protect = $item[4]
    $addr = "$col$r"

    if ($kind -eq 'EMPTY') {
        $ws.Range($addr).ClearContents()
    } elseif ($kind -eq 'NUM') {
        $ws.Range($addr).Value = [double]$val
    } else {
        if ($protect -eq 1) {
            $ws.Range($addr).Value = "'" + $val
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
